# Update "想去人数" (want-to-go count) figures in column F across sheets,
# matching the gh-pages data refresh captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Row = 3;  Value = 986 },
    @{ Sheet = "展览";     Row = 9;  Value = 20 },
    @{ Sheet = "展览";     Row = 10; Value = 389 },
    @{ Sheet = "展览";     Row = 12; Value = 66 },
    @{ Sheet = "展览";     Row = 15; Value = 1954 },
    @{ Sheet = "展览";     Row = 16; Value = 455 },
    @{ Sheet = "展览";     Row = 17; Value = 6645 },
    @{ Sheet = "展览";     Row = 20; Value = 48 },
    @{ Sheet = "展览";     Row = 22; Value = 14 },
    @{ Sheet = "展览";     Row = 23; Value = 204 },
    @{ Sheet = "展览";     Row = 24; Value = 136 },

    @{ Sheet = "演出";     Row = 14; Value = 8 },
    @{ Sheet = "演出";     Row = 18; Value = 28 },
    @{ Sheet = "演出";     Row = 19; Value = 17 },

    @{ Sheet = "本地生活"; Row = 2;  Value = 5448 },
    @{ Sheet = "本地生活"; Row = 4;  Value = 371 },

    @{ Sheet = "全部类型"; Row = 3;  Value = 5448 },
    @{ Sheet = "全部类型"; Row = 5;  Value = 371 },
    @{ Sheet = "全部类型"; Row = 11; Value = 986 },
    @{ Sheet = "全部类型"; Row = 20; Value = 20 },
    @{ Sheet = "全部类型"; Row = 21; Value = 389 },
    @{ Sheet = "全部类型"; Row = 24; Value = 66 },
    @{ Sheet = "全部类型"; Row = 29; Value = 1954 },
    @{ Sheet = "全部类型"; Row = 30; Value = 455 },
    @{ Sheet = "全部类型"; Row = 31; Value = 6645 },
    @{ Sheet = "全部类型"; Row = 35; Value = 48 },
    @{ Sheet = "全部类型"; Row = 36; Value = 85 },
    @{ Sheet = "全部类型"; Row = 37; Value = 8 },
    @{ Sheet = "全部类型"; Row = 38; Value = 14 },
    @{ Sheet = "全部类型"; Row = 39; Value = 204 },
    @{ Sheet = "全部类型"; Row = 41; Value = 136 },
    @{ Sheet = "全部类型"; Row = 44; Value = 28 },
    @{ Sheet = "全部类型"; Row = 45; Value = 17 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
}
